$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the existing header format (bold, centered, bordered)
# from H1 onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-27: column I = I0, column J = IF
$values = @{
    2  = @(6, 7)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(8, 8)
    6  = @(6, 6)
    7  = @(6, 6)
    8  = @(9, 9)
    9  = @(5, 6)
    10 = @(5, 5)
    11 = @(5, 5)
    12 = @(5, 6)
    13 = @(5, 6)
    14 = @(5, 5)
    15 = @(8, 8)
    16 = @(6, 6)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(6, 7)
    20 = @(6, 6)
    21 = @(7, 7)
    22 = @(6, 6)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(7, 8)
    27 = @(7, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$ws.Range("A1").Select()
